# Added high current 3.3v supply
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "high current 3.3v supply" component block in column P
$ws.Range("P12").Value = "Magic Two Transistors"
$ws.Range("P13").Value = "https://www.digikey.com/en/products/detail/rohm-semiconductor/EMH61T2R/5721184"

$ws.Range("P18").Value = "usb controller powersupply"

$ws.Range("P19").Value = "NCP114ASN330T1G"
$ws.Range("P19").Font.Size = 17
$ws.Range("P19").Font.Color = 0x222222
$ws.Range("P19").Font.Name = "Arial"
$ws.Rows(19).RowHeight = 21.75

$ws.Range("P21").Value = "Usb controler"

$ws.Range("O21").Select()
